# Atualiza o link dos fontes de "003%2320180903" para "003%2320180905"
# na ultima linha de texto "Fontes dos exemplos disponiveis em: ..." do slide
# "Material de apoio", quebrando o run do hyperlink em tres runs, exatamente
# como o PowerPoint faz quando o texto de um hyperlink e editado no meio.

$oldFragment = "003%2320180903"
$newMiddle   = "/003%2320180905/"
$newFirst    = "https://github.com/plsqlcamp/Meetup/tree/master"
$newLast     = "Source"

$p = $ppt.ActivePresentation

$targetShape = $null
$targetSlide = $null
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $s.Shapes.Count; $shi++) {
        $sh = $s.Shapes.Item($shi)
        if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
            if ($sh.TextFrame.TextRange.Text -like "*$oldFragment*") {
                $targetShape = $sh
                $targetSlide = $s
            }
        }
    }
}

$tf = $targetShape.TextFrame
$tr = $tf.TextRange

$targetParagraph = $null
for ($i = 1; $i -le $tr.Paragraphs().Count; $i++) {
    $para = $tr.Paragraphs($i)
    if ($para.Text -like "*$oldFragment*") {
        $targetParagraph = $para
    }
}

$targetRunIdx = -1
for ($i = 1; $i -le $targetParagraph.Runs().Count; $i++) {
    $run = $targetParagraph.Runs($i)
    if ($run.Text -like "*$oldFragment*") {
        $targetRunIdx = $i
    }
}
$targetRun = $targetParagraph.Runs($targetRunIdx)

# O texto do ultimo run de um paragrafo inclui a marca de fim de paragrafo,
# entao usamos o comprimento real da URL (sem essa marca) para recortar via
# Characters, preservando a formatacao/hyperlink original do run.
$fullLen = $targetRun.Text.Length
$urlText = $targetRun.Text.Substring(0, $fullLen - 1)
if ($urlText.Substring($urlText.Length - 1, 1) -eq [char]13) {
    $urlText = $urlText.Substring(0, $urlText.Length - 1)
}
$urlLen = $urlText.Length

# 1) Primeiro run: encurta o texto existente via Characters (mantem
#    lang/sz/dirty/hyperlink originais do run).
$urlRange = $tr.Characters($targetRun.Start, $urlLen)
$urlRange.Text = $newFirst

# Reconsultamos o run (por indice) em vez de encadear a partir do objeto
# Characters(...) devolvido acima: TextRange.InsertAfter so reproduz o
# comportamento de "continuar digitando dentro do hyperlink" quando chamado
# sobre o proprio objeto Run.
$firstRun = $targetParagraph.Runs($targetRunIdx)

# 2) Segundo run: insere a nova parte da data logo apos o primeiro
$secondRun = $firstRun.InsertAfter($newMiddle)

# 3) Terceiro run: insere "Source" logo apos o segundo
$thirdRun = $secondRun.InsertAfter($newLast)
